$wb = $excel.ActiveWorkbook

# --- Exam Dashboard: update remark text now that the date check passed ---
$examWs = $wb.Worksheets.Item("Exam Dashboard")
$examWs.Range("E3").Value = "date is valid"
$examWs.Range("E4").Value = "date is valid"

# Comments column no longer needs to be as wide with the shorter message
$examWs.Columns.Item(5).ColumnWidth = 14.1666666666667

# --- Recolor the bold title/header fonts to white on every sheet ---
foreach ($sheet in $wb.Worksheets) {
    # Title cell (row 1): drop the old 14pt size and switch the bold font to white
    $titleFont = $sheet.Range("A1").Font
    $titleFont.Size = 11
    $titleFont.Color = 16777215

    # Header row (row 2): keep bold, switch the font color to white.
    # Limit to the sheet's own used columns so we don't touch untouched cells.
    $lastCol = $sheet.UsedRange.Columns.Count
    $headerRange = $sheet.Range($sheet.Cells.Item(2, 1), $sheet.Cells.Item(2, $lastCol))
    $headerRange.Font.Color = 16777215
}
